$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Rename existing table/header columns to the new (English) names ---
# (writing the header cell value keeps the ListObject column name in sync)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "User"

# --- Add two new columns to the query table: Group, Policy ---
$colGroup  = $tbl.ListColumns.Add()
$colPolicy = $tbl.ListColumns.Add()
$ws.Range("E1").Value = "Group"
$ws.Range("F1").Value = "Policy"

# Give the two new columns a sensible width like the rest of the sheet
$ws.Columns.Item(5).ColumnWidth = 16.29
$ws.Columns.Item(6).ColumnWidth = 17.29

# --- Populate Group / Policy for every data row ---
# Rows 2-12  -> CloudAdmins
# Rows 13-18 -> CloudDBAs
# Rows 19-31 -> CloudReadOnly
# All rows get the same enforced-MFA policy
for ($r = 2; $r -le 12; $r++) {
  $ws.Cells.Item($r, 5).Value = "CloudAdmins"
  $ws.Cells.Item($r, 6).Value = "EnforceMFAPolicy"
}
for ($r = 13; $r -le 18; $r++) {
  $ws.Cells.Item($r, 5).Value = "CloudDBAs"
  $ws.Cells.Item($r, 6).Value = "EnforceMFAPolicy"
}
for ($r = 19; $r -le 31; $r++) {
  $ws.Cells.Item($r, 5).Value = "CloudReadOnly"
  $ws.Cells.Item($r, 6).Value = "EnforceMFAPolicy"
}

# --- Move the selection to reflect where the author ended up editing ---
$ws.Range("A31").Select()
